# Update for release to deploy 0.1.1
#
# Changes:
#  1. Rename the second worksheet from "Include from NMDP Disease Sta" to "Include #0"
#  2. On the "Metadata" sheet:
#       - Insert a new "Jurisdiction" row (empty value) right after the "Contact" row,
#         pushing "Description", "Purpose", "Copyright" and "Immutable" down by one row.
#       - Update the "Version" value from 0.1.0 to 0.1.1
#       - Update the "Date" value to the new timestamp

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include..." sheet (the non-Metadata sheet) ---
$wsInclude = $wb.Worksheets.Item(2)
if ($wsInclude.Name -eq "Metadata") {
  $wsInclude = $wb.Worksheets.Item(1)
}
$wsInclude.Name = "Include #0"

# --- 2. Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Insert a blank row above the current row 11 ("Description"), shifting rows down.
$ws.Rows.Item(11).Insert()

# Copy the formatting of the (now shifted) "Description" row onto the new blank row
# so the new row matches the rest of the table's style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new "Jurisdiction" row.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Bump the Version property value.
$ws.Range("B3").Value = "0.1.1"

# Update the Date property value.
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"
